# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect the latest generated counts (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Mapping of row number -> new value for column F
$updates = @{
    2  = 1111
    3  = 822
    4  = 275
    5  = 46
    6  = 1108
    8  = 2043
    9  = 7620
    10 = 905
    11 = 425
    12 = 352
    14 = 405
    15 = 154
    16 = 7772
    17 = 308
    18 = 1346
    19 = 152
    21 = 228
    22 = 154
    23 = 310
    24 = 146
    25 = 164
    26 = 19
    27 = 107
    28 = 23
    29 = 413
    30 = 620
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
